# Sprint 2 Acceptance tests: record results for the "Acceptance test" task
# row (Assigned to / Dates / Taken) and the story's "Total Hours" Taken cell
# in each of the five Task tables touched by this commit.
#
# Table layout (columns): Task ID | Task Description | Assigned to | Dates |
# Estimate | Taken
# The "Acceptance test" task is always the second-to-last row of its table,
# immediately followed by the "Story Points / Total Hours" summary row.

$d = $word.ActiveDocument

# 1-based Word Tables.Item() index -> [Owner, Date, AcceptanceTestTaken, TotalHoursTaken]
$updates = @{
    13 = @("Victoria Tobin", "20/09", "0.5", "3.5")   # Story A-1: All Orders
    14 = @("Victoria Tobin", "20/09", "1",   "2")     # Story A-2: Assigned Courier
    15 = @("Victoria Tobin", "20/09", "0.5", "3.5")   # Story C-3: Order Status
    16 = @("Victoria Tobin", "20/09", "0.5", "3")     # Story C-5: Order Form
    17 = @("Victoria Tobin", "20/09", "0.5", "1")     # Story C-4: Pick-Up Request
}

foreach ($tableIndex in $updates.Keys) {
    $values = $updates[$tableIndex]
    $owner = $values[0]
    $date = $values[1]
    $testTaken = $values[2]
    $totalTaken = $values[3]

    $table = $d.Tables.Item($tableIndex)
    $rowCount = $table.Rows.Count

    $acceptanceRow = $table.Rows.Item($rowCount - 1)
    $totalsRow = $table.Rows.Item($rowCount)

    $acceptanceRow.Cells.Item(3).Range.Text = $owner
    $acceptanceRow.Cells.Item(4).Range.Text = $date
    $acceptanceRow.Cells.Item(6).Range.Text = $testTaken

    $totalsRow.Cells.Item(6).Range.Text = $totalTaken
}

Write-Output "done"
